# Apply updated B column values (rows 2-101) on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    -121.1398329593187,
    -119.4150302107065,
    -117.6902274620943,
    -115.9654247134822,
    -114.24062196487,
    -112.5158192162579,
    -110.7910164676457,
    -109.0662137190335,
    -107.3414109704215,
    -105.6166082218093,
    -103.8918054731971,
    -102.167002724585,
    -100.4421999759728,
    -98.71739722736072,
    -96.99259447874852,
    -95.26779173013631,
    -93.54298898152422,
    -91.81818623291201,
    -90.09338348429992,
    -88.36858073568771,
    -86.64377798707551,
    -84.91897523846342,
    -83.19417248985121,
    -81.46936974123912,
    -79.74456699262691,
    -78.01976424401482,
    -76.29496149540262,
    -74.57015874679041,
    -72.84535599817832,
    -71.12055324956623,
    -69.39575050095402,
    -67.67094775234193,
    -65.94614500372973,
    -64.22134225511752,
    -62.49653950650543,
    -60.77173675789322,
    -59.04693400928113,
    -57.32213126066893,
    -55.59732851205683,
    -53.87252576344463,
    -52.14772301483242,
    -50.42292026622033,
    -48.69811751760813,
    -46.97331476899603,
    -45.24851202038383,
    -43.52370927177174,
    -41.79890652315953,
    -40.07410377454744,
    -38.34930102593523,
    -36.62449827732314,
    -34.89969552871094,
    -33.17489278009873,
    -31.45009003148664,
    -29.72528728287443,
    -28.00048453426234,
    -26.27568178565014,
    -24.55087903703804,
    -22.82607628842584,
    -21.10127353981375,
    -19.37647079120154,
    -17.65166804258934,
    -15.92686529397724,
    -14.20206254536504,
    -12.47725979675295,
    -10.75245704814074,
    -9.027654299528649,
    -7.302851550916444,
    -5.578048802304295,
    -3.853246053692146,
    -2.128443305079941,
    -0.4036405564677921,
    1.321162192144357,
    3.045964940756448,
    4.770767689368654,
    6.495570437980746,
    8.220373186592951,
    9.9451759352051,
    11.66997868381725,
    13.3947814324294,
    15.11958418104155,
    16.84438692965369,
    18.56918967826584,
    20.29399242687805,
    22.0187951754902,
    23.74359792410235,
    25.46840067271444,
    27.19320342132664,
    28.91800616993885,
    30.64280891855094,
    32.36761166716315,
    34.09241441577524,
    35.81721716438744,
    37.54201991299954,
    39.26682266161174,
    40.99162541022389,
    42.71642815883604,
    44.44123090744813,
    46.16603365606034,
    47.89083640467248,
    49.61563915328463
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}
